# Apply the "Fixes up until clusterrolebinding" edit:
# - "clusters" sheet: replace the two sample rows with a single testcluster row
# - "components" sheet: no data change, just selection/active-cell bookkeeping
# - "environments" sheet: unchanged
# - workbook view: window geometry changed, first sheet (clusters) becomes active

$wb = $excel.ActiveWorkbook

$wsClusters = $wb.Worksheets.Item("clusters")
$wsComponents = $wb.Worksheets.Item("components")
$wsEnvironments = $wb.Worksheets.Item("environments")

# --- clusters sheet: collapse to a single "testcluster" row ---
$wsClusters.Rows.Item(2).Delete()

$wsClusters.Range("A1").Value = "testcluster"
$wsClusters.Range("B1").Value = "Een testcluster"
$wsClusters.Range("C1").Value = "Een domein voor het testcluster"
$wsClusters.Range("D1").Value = "a.conduction.nl"

# --- selection / active sheet bookkeeping ---
$wsComponents.Activate() | Out-Null
$wsComponents.Range("D8").Select() | Out-Null

$wsClusters.Activate() | Out-Null
$wsClusters.Range("A1:D1").Select() | Out-Null
